# Sync attendance_reports - apply the "Recorded By" reordering, updated
# statistics, and the newly-recorded PARASITOLOGY session (6) for Year 2 / C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2..5 (ANATOMY): "Recorded By" lists reordered ---
$ws.Range("G2").Value = "System, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value = "System, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 6: Recorded Sessions metric 23 -> 24 ---
$ws.Range("L6").Value = 24

# --- Row 7 (BIOCHEMISTRY LAB/CBL #1): reorder + Missing Sessions 3 -> 2 ---
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg"
$ws.Range("L7").Value = 2

# --- Row 9 (HISTOLOGY #1): reorder + Coverage % 79.3% -> 82.8% ---
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "82.8%"
$ws.Range("L8").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# --- Row 10: Average Attendance % 27.1% -> 27.4% ---
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "27.4%"
$ws.Range("L8").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- Row 12 (MICROBIOLOGY #1): reorder ---
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg"

# --- Row 15 (Group Statistics summary row): Recorded/Missing/Coverage/Avg ---
$ws.Range("O15").Value = 24
$ws.Range("P15").Value = 2
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "82.8%"
$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "27.4%"
$ws.Range("Q15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# --- Row 17 (PARASITOLOGY #5): reorder ---
$ws.Range("G17").Value = "esraa.sami@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

# --- Row 18 (PARASITOLOGY #6): session newly recorded ---
# Copy the "Recorded" (green) formatting from row 17 onto row 18 first,
# then fill in the recorder, student count and status.
$ws.Range("A17:I17").Copy()
$ws.Range("A18:I18").PasteSpecial(-4122)
$ws.Range("G18").Value = "afnan.fares@med.asu.edu.eg"
$ws.Range("H18").Value = "83/251"
$ws.Range("I18").Value = "Recorded"

# --- Row 24 (PATHOLOGY LAB/MUSEUM #2): reorder ---
$ws.Range("G24").Value = "Sarah.Mahdy@med.asu.edu.eg, youstina.gamil@med.asu.edu.eg"

# --- Row 27 (PHARMACOLOGY #2): reorder ---
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 30 (PHYSIOLOGY #3): reorder ---
$ws.Range("G30").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
